$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.279.59"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "3.140.25"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.135.77"
$ws.Range("E8").Value = "  +3.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.486"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.124"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "3.659.13"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "67.246.81"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "3.139.52"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.78%  "
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "2.863.98"
$ws.Range("E45").Value = "  +5.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "389.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("E51").Value = "  +0.18%  "
